$d = $word.ActiveDocument

# Locate the "DWL^PC = 0 < DWL^M" run. In the target revision this text is
# split into three runs: the original lead-in text (now with a trailing
# space and explicit xml:space="preserve"), then "DWL^M" and "Wha" each in
# their own run, wrapped by a proofErr spellStart/spellEnd pair (as Word's
# spell checker does for an out-of-dictionary word like "DWL^MWha").
$rng = $d.Content
$found = $rng.Find.Execute("DWL^^PC = 0 < DWL^^M")
if (-not $found) {
    throw "Could not find target text 'DWL^PC = 0 < DWL^M'"
}

$target = $d.Range($rng.Start, $rng.End)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
       '<w:body>' +
       '<w:p>' +
       '<w:r w:rsidRPr="0018554B">' +
       '<w:rPr>' +
       '<w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>' +
       '<w:kern w:val="0"/>' +
       '<w:sz w:val="21"/>' +
       '<w:szCs w:val="21"/>' +
       '<w14:ligatures w14:val="none"/>' +
       '</w:rPr>' +
       '<w:t xml:space="preserve">DWL^PC = 0 &lt; </w:t>' +
       '</w:r>' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r>' +
       '<w:rPr>' +
       '<w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>' +
       '<w:kern w:val="0"/>' +
       '<w:sz w:val="21"/>' +
       '<w:szCs w:val="21"/>' +
       '<w14:ligatures w14:val="none"/>' +
       '</w:rPr>' +
       '<w:t>DWL^M</w:t>' +
       '</w:r>' +
       '<w:r>' +
       '<w:rPr>' +
       '<w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/>' +
       '<w:kern w:val="0"/>' +
       '<w:sz w:val="21"/>' +
       '<w:szCs w:val="21"/>' +
       '<w14:ligatures w14:val="none"/>' +
       '</w:rPr>' +
       '<w:t>Wha</w:t>' +
       '</w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '</w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData>' +
       '</pkg:part>' +
       '</pkg:package>'

$target.InsertXML($xml)
